$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Angpt1"
$ws.Cells.Item(2, 3).Value = "Tek"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.07171233333333334
$ws.Cells.Item(2, 8).Value = 0.215137
$ws.Cells.Item(2, 9).Value = 0.004947717871829782
$ws.Cells.Item(2, 10).Value = 0.004947717871829783
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 44.69746666666666
$ws.Cells.Item(2, 14).Value = 134.0924
$ws.Cells.Item(2, 15).Value = 0.6823972194925493
$ws.Cells.Item(2, 16).Value = 0.6823972194925493
$ws.Cells.Item(2, 17).Value = 3.205359628755555
$ws.Cells.Item(2, 18).Value = 28.8482366588
$ws.Cells.Item(2, 19).Value = 0.003376308918570237
$ws.Cells.Item(2, 20).Value = 0.003376308918570237

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Angpt1"
$ws.Cells.Item(3, 3).Value = "Tek"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.07171233333333334
$ws.Cells.Item(3, 8).Value = 0.215137
$ws.Cells.Item(3, 9).Value = 0.004947717871829782
$ws.Cells.Item(3, 10).Value = 0.004947717871829783
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 18.63243533333333
$ws.Cells.Item(3, 14).Value = 55.897306
$ws.Cells.Item(3, 15).Value = 0.2844618053784121
$ws.Cells.Item(3, 16).Value = 0.2844618053784121
$ws.Cells.Item(3, 17).Value = 1.336175413435778
$ws.Cells.Item(3, 18).Value = 12.025578720922
$ws.Cells.Item(3, 19).Value = 0.001407436758323735
$ws.Cells.Item(3, 20).Value = 0.001407436758323735

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Angpt1"
$ws.Cells.Item(4, 3).Value = "Tek"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.07171233333333334
$ws.Cells.Item(4, 8).Value = 0.215137
$ws.Cells.Item(4, 9).Value = 0.004947717871829782
$ws.Cells.Item(4, 10).Value = 0.004947717871829783
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.170755666666667
$ws.Cells.Item(4, 14).Value = 6.512267
$ws.Cells.Item(4, 15).Value = 0.03314097512903853
$ws.Cells.Item(4, 16).Value = 0.03314097512903853
$ws.Cells.Item(4, 17).Value = 0.1556699539532223
$ws.Cells.Item(4, 18).Value = 1.401029585579
$ws.Cells.Item(4, 19).Value = 0.0001639721949358103
$ws.Cells.Item(4, 20).Value = 0.0001639721949358102

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Angpt1"
$ws.Cells.Item(5, 3).Value = "Tek"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 11.05178533333333
$ws.Cells.Item(5, 8).Value = 33.155356
$ws.Cells.Item(5, 9).Value = 0.7625064374239614
$ws.Cells.Item(5, 10).Value = 0.7625064374239615
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 44.69746666666666
$ws.Cells.Item(5, 14).Value = 134.0924
$ws.Cells.Item(5, 15).Value = 0.6823972194925493
$ws.Cells.Item(5, 16).Value = 0.6823972194925493
$ws.Cells.Item(5, 17).Value = 493.9868065438222
$ws.Cells.Item(5, 18).Value = 4445.881258894399
$ws.Cells.Item(5, 19).Value = 0.5203322727432808
$ws.Cells.Item(5, 20).Value = 0.5203322727432809

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Angpt1"
$ws.Cells.Item(6, 3).Value = "Tek"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.05178533333333
$ws.Cells.Item(6, 8).Value = 33.155356
$ws.Cells.Item(6, 9).Value = 0.7625064374239614
$ws.Cells.Item(6, 10).Value = 0.7625064374239615
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 18.63243533333333
$ws.Cells.Item(6, 14).Value = 55.897306
$ws.Cells.Item(6, 15).Value = 0.2844618053784121
$ws.Cells.Item(6, 16).Value = 0.2844618053784121
$ws.Cells.Item(6, 17).Value = 205.9216755412151
$ws.Cells.Item(6, 18).Value = 1853.295079870936
$ws.Cells.Item(6, 19).Value = 0.2169039578022813
$ws.Cells.Item(6, 20).Value = 0.2169039578022813

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Angpt1"
$ws.Cells.Item(7, 3).Value = "Tek"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.05178533333333
$ws.Cells.Item(7, 8).Value = 33.155356
$ws.Cells.Item(7, 9).Value = 0.7625064374239614
$ws.Cells.Item(7, 10).Value = 0.7625064374239615
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.170755666666667
$ws.Cells.Item(7, 14).Value = 6.512267
$ws.Cells.Item(7, 15).Value = 0.03314097512903853
$ws.Cells.Item(7, 16).Value = 0.03314097512903853
$ws.Cells.Item(7, 17).Value = 23.99072563911689
$ws.Cells.Item(7, 18).Value = 215.916530752052
$ws.Cells.Item(7, 19).Value = 0.02527020687839928
$ws.Cells.Item(7, 20).Value = 0.02527020687839928

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Angpt1"
$ws.Cells.Item(8, 3).Value = "Tek"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.370524666666666
$ws.Cells.Item(8, 8).Value = 10.111574
$ws.Cells.Item(8, 9).Value = 0.2325458447042087
$ws.Cells.Item(8, 10).Value = 0.2325458447042088
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 44.69746666666666
$ws.Cells.Item(8, 14).Value = 134.0924
$ws.Cells.Item(8, 15).Value = 0.6823972194925493
$ws.Cells.Item(8, 16).Value = 0.6823972194925493
$ws.Cells.Item(8, 17).Value = 150.6539139375111
$ws.Cells.Item(8, 18).Value = 1355.8852254376
$ws.Cells.Item(8, 19).Value = 0.1586886378306982
$ws.Cells.Item(8, 20).Value = 0.1586886378306982

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Angpt1"
$ws.Cells.Item(9, 3).Value = "Tek"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.370524666666666
$ws.Cells.Item(9, 8).Value = 10.111574
$ws.Cells.Item(9, 9).Value = 0.2325458447042087
$ws.Cells.Item(9, 10).Value = 0.2325458447042088
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 18.63243533333333
$ws.Cells.Item(9, 14).Value = 55.897306
$ws.Cells.Item(9, 15).Value = 0.2844618053784121
$ws.Cells.Item(9, 16).Value = 0.2844618053784121
$ws.Cells.Item(9, 17).Value = 62.80108289107155
$ws.Cells.Item(9, 18).Value = 565.209746019644
$ws.Cells.Item(9, 19).Value = 0.06615041081780708
$ws.Cells.Item(9, 20).Value = 0.06615041081780709

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Angpt1"
$ws.Cells.Item(10, 3).Value = "Tek"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.370524666666666
$ws.Cells.Item(10, 8).Value = 10.111574
$ws.Cells.Item(10, 9).Value = 0.2325458447042087
$ws.Cells.Item(10, 10).Value = 0.2325458447042088
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.170755666666667
$ws.Cells.Item(10, 14).Value = 6.512267
$ws.Cells.Item(10, 15).Value = 0.03314097512903853
$ws.Cells.Item(10, 16).Value = 0.03314097512903853
$ws.Cells.Item(10, 17).Value = 7.316585519806445
$ws.Cells.Item(10, 18).Value = 65.849269678258
$ws.Cells.Item(10, 19).Value = 0.007706796055703439
$ws.Cells.Item(10, 20).Value = 0.007706796055703439
